# Update "想去人数" (want-to-go count) figures that changed between scrapes.
#
# 展览 (Exhibitions) sheet:
#   F6  (合肥·环形宇宙动漫游戏嘉年华) : 2957 -> 3030
#   F13 (合肥·首届运动番only)        : 11   -> 12
#
# 全部类型 (All types) sheet - same two events, shifted one row down:
#   F7  (合肥·环形宇宙动漫游戏嘉年华) : 2957 -> 3030
#   F14 (合肥·首届运动番only)        : 11   -> 12

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F6").Value = 3030
$wsExhibition.Range("F13").Value = 12

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F7").Value = 3030
$wsAllTypes.Range("F14").Value = 12
